$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New user rows to append below the existing data (rows 2-5 are populated already).
$rows = @(
    @{ A = "mohamed.hussein.eladwy"; B = "Mohamed Eladwy"; C = "01123456789"; D = "mh8579007@gmail.com";        E = "Nlsznvw123*+" },
    @{ A = "mohamed.eladwy";         B = "Mohamed Eladwy"; C = "01233456789"; D = "mohamedeljoker309@gmail.com"; E = "Nlsznvw123#" },
    @{ A = "yusuf.bdr";              B = "Yusuf Elsayed";  C = "01223456789"; D = "yusuf.bdr@gmail.com";         E = "BfhfuYwi123@" }
)

$startRow = 6
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $false

    # Columns C and E carry a column-level style (numFmtId 49) that should
    # not bleed into these new rows, so reset them back to the default style.
    $ws.Cells.Item($r, 3).Style = "Normal"
    $ws.Cells.Item($r, 5).Style = "Normal"
}

# Set column A and B widths individually (mirrors the col split seen in the diff).
# 24.17 round-trips to a stored width of 25, matching the original formatting.
$ws.Columns.Item(1).ColumnWidth = 24.17
$ws.Columns.Item(2).ColumnWidth = 24.17
